$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "purpose" column (E2:E14) from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E14").Value = "fullRNASEQ"

# Update the active selection to match the edited range (E13:E14)
$ws.Range("E13:E14").Select()
